$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) New styled cells that reuse the existing "fillId5 + border" look
#    (the same format already used by E19 in the original workbook).
#    Grab that format once and stamp it onto the other cells that need it.
# ------------------------------------------------------------------
$ws.Range("E19").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E14").PasteSpecial(-4122)

# E14 additionally gets a left-aligned version of that style
$ws.Range("E14").HorizontalAlignment = -4131

# ------------------------------------------------------------------
# 2) New CLOSEALERT step (row 14): Keyword + resulting alert message
# ------------------------------------------------------------------
$ws.Range("B14").Value = "CLOSEALERT"
$ws.Range("E14").Value = "You Have Succesfully Logged Out!!"

# ------------------------------------------------------------------
# 3) VERIFYTITLE "Handling 404": E17 had a dead short URL with its own
#    hyperlink + special styling; make it a plain value matching the
#    correct long URL (same as E3/E8/E23), no hyperlink.
# ------------------------------------------------------------------
$ws.Range("E17").Hyperlinks.Delete()
$ws.Range("A17").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = "http://demo.guru99.com/V4/"

# ------------------------------------------------------------------
# 4) E20 no longer carries any formatting/value - fully clear it.
# ------------------------------------------------------------------
$ws.Range("E20").ClearContents()
$ws.Range("E20").ClearFormats()

# ------------------------------------------------------------------
# 5) Remove the obsolete "Check Lenta" test row (old row 25: lenta/css)
#    by deleting the whole row - everything below shifts up.
# ------------------------------------------------------------------
$ws.Rows("25:25").Delete()

# ------------------------------------------------------------------
# 6) View tidy-up to match the saved workbook view state.
# ------------------------------------------------------------------
$ws.Range("E26").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.Zoom = 100

Write-Host "edit applied"
